$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = "r27"
$ws.Range("B28").Value = "paciente mayor de edad"
$ws.Range("C28").Value = "!is.na(fechaid) & !is.na(interview) & (interval(fechaid, interview) / years(1)) >= 18"
